$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H mirrors the existing "Expansion Cost" block in column G but
# for a new "Expansion Trafo" concept. Copy each source cell (G) onto the
# new cell (H) first so the new column picks up the same formatting/style,
# then overwrite the cell values that differ from column G.

$ws.Range("G1").Copy($ws.Range("H1"))

$ws.Range("G3").Copy($ws.Range("H3"))
$ws.Range("H3").Value = "Expansion  Trafo"

$ws.Range("G4").Copy($ws.Range("H4"))
$ws.Range("H4").Value = "pExpTrafoP"

$ws.Range("G5").Copy($ws.Range("H5"))
$ws.Range("H5").Value = "Additional Link Power with one additional Trafo"

$ws.Range("G6").Copy($ws.Range("H6"))

$ws.Range("G7").Copy($ws.Range("H7"))
$ws.Range("H7").Value = "[MW]"

$ws.Range("G8").Copy($ws.Range("H8"))
$ws.Range("H8").Value = 50

$ws.Range("G9").Copy($ws.Range("H9"))
$ws.Range("H9").Value = 50

$ws.Columns.Item(8).ColumnWidth = 23

$excel.Goto($ws.Range("H4"))

Write-Output "done"
